# Updated cryptos list on Tue Sep 19 07:42:09 UTC 2023 with GitHub Actions
# Refreshes Price (D) / Volume(1h) (E) figures for most rows, and shifts
# rows 48-51 down by one to make room for the newly-listed BabyDogeCoin
# (the former last row, USDD, drops off the bottom of the table).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.891.90"
$ws.Range("E2").Value = "  +0.13%  "

$ws.Range("D3").Value = "1.639.26"
$ws.Range("E3").Value = "  -0.20%  "

$ws.Range("E4").Value = "  -0.47%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "217.13"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.46%  "

$ws.Range("E6").Value = "  +1.81%  "

$ws.Range("E7").Value = "  -0.45%  "

$ws.Range("E8").Value = "  +0.89%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.0623"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +0.14%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.88"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +3.28%  "

$ws.Range("E11").Value = "  +0.11%  "

$ws.Range("E12").Value = "  -0.21%  "

$ws.Range("D13").Value = "1.656.40"
$ws.Range("E13").Value = "  +0.79%  "

$ws.Range("E14").Value = "  -1.14%  "

$ws.Range("E15").Value = "  +0.67%  "

$ws.Range("E16").Value = "  +2.70%  "

$ws.Range("D17").Value = "26.882.01"
$ws.Range("E17").Value = "  +0.00%  "

$ws.Range("E18").Value = "  +0.03%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "218.23"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +0.93%  "

$ws.Range("E20").Value = "  -0.37%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.70"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +1.82%  "

$ws.Range("E22").Value = "  +0.67%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.42"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +2.58%  "

$ws.Range("E24").Value = "  -0.34%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "147.05"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -0.37%  "

$ws.Range("E26").Value = "  -0.53%  "

$ws.Range("E27").Value = "  +0.48%  "

$ws.Range("E28").Value = "  +1.14%  "

$ws.Range("E29").Value = "  -0.08%  "

$ws.Range("E30").Value = "  -1.21%  "

$ws.Range("E31").Value = "  -1.08%  "

$ws.Range("E32").Value = "  -1.22%  "

$ws.Range("E33").Value = "  +0.00%  "

$ws.Range("E34").Value = "  +0.85%  "

$ws.Range("D35").Value = "1.264.95"
$ws.Range("E35").Value = "  -1.46%  "

$ws.Range("E36").Value = "  +0.05%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0178"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +2.48%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.838"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +1.84%  "

$ws.Range("E39").Value = "  +0.16%  "

$ws.Range("E40").Value = "  -0.45%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.811"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +0.53%  "

$ws.Range("E42").Value = "  +0.63%  "

$ws.Range("D43").Value = "1.779.57"
$ws.Range("E43").Value = "  -0.17%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "62.30"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +1.86%  "

$ws.Range("E45").Value = "  +0.37%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "92.01"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -0.63%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.61"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +0.68%  "

$ws.Range("B48").Value = "BabyDogeCoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D48").Value = "0.0₆0105"
$ws.Range("E48").Value = "  +10.36%  "

$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0512"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -0.91%  "

$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.68"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +1.14%  "

$ws.Range("B51").Value = "Algorand"
$ws.Range("C51").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0962"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -0.78%  "
